# "added final presentation v2" - removes the last slide (id 296, the
# "Screenshots" slide) and its notes page from the deck.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item($p.Slides.Count)
$s.Delete()
